$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to Text format so that
# numeric-looking strings (e.g. "61.320.60", "1.00", "43.00") are preserved
# exactly as text instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Apply the updated crypto price / volume figures cell by cell
$ws.Cells.Item(2, 4).Value = '61.320.60'
$ws.Cells.Item(2, 5).Value = '  +2.85%  '
$ws.Cells.Item(3, 4).Value = '3.401.38'
$ws.Cells.Item(3, 5).Value = '  +4.94%  '
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.55%  '
$ws.Cells.Item(5, 4).Value = '405.79'
$ws.Cells.Item(5, 5).Value = '  +0.76%  '
$ws.Cells.Item(6, 4).Value = '128.76'
$ws.Cells.Item(6, 5).Value = '  +17.94%  '
$ws.Cells.Item(7, 4).Value = '0.607'
$ws.Cells.Item(7, 5).Value = '  +9.61%  '
$ws.Cells.Item(8, 4).Value = '3.393.57'
$ws.Cells.Item(8, 5).Value = '  +4.74%  '
$ws.Cells.Item(9, 5).Value = '  +0.09%  '
$ws.Cells.Item(10, 4).Value = '0.673'
$ws.Cells.Item(10, 5).Value = '  +11.88%  '
$ws.Cells.Item(11, 4).Value = '0.127'
$ws.Cells.Item(11, 5).Value = '  +22.40%  '
$ws.Cells.Item(12, 4).Value = '42.29'
$ws.Cells.Item(12, 5).Value = '  +12.35%  '
$ws.Cells.Item(13, 5).Value = '  -0.36%  '
$ws.Cells.Item(14, 4).Value = '3.950.71'
$ws.Cells.Item(14, 5).Value = '  +4.54%  '
$ws.Cells.Item(15, 4).Value = '8.55'
$ws.Cells.Item(15, 5).Value = '  +7.45%  '
$ws.Cells.Item(16, 4).Value = '19.66'
$ws.Cells.Item(16, 5).Value = '  +6.11%  '
$ws.Cells.Item(17, 4).Value = '3.398.76'
$ws.Cells.Item(17, 5).Value = '  +3.34%  '
$ws.Cells.Item(18, 4).Value = '11.63'
$ws.Cells.Item(18, 5).Value = '  +13.35%  '
$ws.Cells.Item(19, 4).Value = '61.394.97'
$ws.Cells.Item(19, 5).Value = '  +2.53%  '
$ws.Cells.Item(20, 5).Value = '  +5.04%  '
$ws.Cells.Item(21, 4).Value = '0.0000135'
$ws.Cells.Item(21, 5).Value = '  +25.86%  '
$ws.Cells.Item(22, 5).Value = '  +3.17%  '
$ws.Cells.Item(23, 4).Value = '82.77'
$ws.Cells.Item(23, 5).Value = '  +14.89%  '
$ws.Cells.Item(24, 4).Value = '13.09'
$ws.Cells.Item(24, 5).Value = '  +9.86%  '
$ws.Cells.Item(25, 4).Value = '307.54'
$ws.Cells.Item(25, 5).Value = '  +6.23%  '
$ws.Cells.Item(26, 4).Value = '3.16'
$ws.Cells.Item(26, 5).Value = '  +5.14%  '
$ws.Cells.Item(27, 4).Value = '8.59'
$ws.Cells.Item(27, 5).Value = '  +17.62%  '
$ws.Cells.Item(28, 4).Value = '4.72'
$ws.Cells.Item(28, 5).Value = '  +5.39%  '
$ws.Cells.Item(29, 4).Value = '29.67'
$ws.Cells.Item(29, 5).Value = '  +6.69%  '
$ws.Cells.Item(30, 5).Value = '  +4.49%  '
$ws.Cells.Item(31, 5).Value = '  +5.32%  '
$ws.Cells.Item(32, 5).Value = '  +9.04%  '
$ws.Cells.Item(33, 4).Value = '11.71'
$ws.Cells.Item(33, 5).Value = '  +7.33%  '
$ws.Cells.Item(34, 4).Value = '43.00'
$ws.Cells.Item(34, 5).Value = '  +12.65%  '
$ws.Cells.Item(35, 5).Value = '  +11.78%  '
$ws.Cells.Item(36, 5).Value = '  +0.26%  '
$ws.Cells.Item(37, 4).Value = '0.0487'
$ws.Cells.Item(37, 5).Value = '  +5.99%  '
$ws.Cells.Item(38, 4).Value = '52.15'
$ws.Cells.Item(38, 5).Value = '  +0.53%  '
$ws.Cells.Item(39, 4).Value = '0.998'
$ws.Cells.Item(39, 5).Value = '  -0.78%  '
$ws.Cells.Item(40, 5).Value = '  +6.28%  '
$ws.Cells.Item(41, 4).Value = '3.00'
$ws.Cells.Item(41, 5).Value = '  +0.64%  '
$ws.Cells.Item(42, 5).Value = '  +7.44%  '
$ws.Cells.Item(43, 4).Value = '1.97'
$ws.Cells.Item(43, 5).Value = '  +8.26%  '
$ws.Cells.Item(44, 4).Value = '136.41'
$ws.Cells.Item(44, 5).Value = '  +1.53%  '
$ws.Cells.Item(45, 2).Value = 'TheGraph'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(45, 4).Value = '0.285'
$ws.Cells.Item(45, 5).Value = '  +5.92%  '
$ws.Cells.Item(46, 2).Value = 'NEARProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(46, 4).Value = '3.94'
$ws.Cells.Item(46, 5).Value = '  +7.81%  '
$ws.Cells.Item(47, 4).Value = '16.95'
$ws.Cells.Item(47, 5).Value = '  +7.68%  '
$ws.Cells.Item(48, 5).Value = '  +2.54%  '
$ws.Cells.Item(49, 4).Value = '21.87'
$ws.Cells.Item(49, 5).Value = '  +8.08%  '
$ws.Cells.Item(50, 4).Value = '2.147.15'
$ws.Cells.Item(50, 5).Value = '  +2.36%  '
$ws.Cells.Item(51, 4).Value = '3.742.09'
$ws.Cells.Item(51, 5).Value = '  +2.47%  '
